$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. TestBase sheet: flip the Regression (column C) flags for the
#    existing Registration/Login rows, and append a new MyAccount row.
# ---------------------------------------------------------------------
$testBase = $wb.Worksheets.Item("TestBase")
$testBase.Range("C2").Value = "YES"
$testBase.Range("C3").Value = "NO"
$testBase.Range("A4").Value = "MyAccount"
$testBase.Range("B4").Value = "NO"
$testBase.Range("C4").Value = "NO"
$testBase.Range("D4").Value = "NO"

# ---------------------------------------------------------------------
# 2. Select the whole Login sheet (mirrors the saved selection state
#    captured in the target workbook for that tab).
# ---------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Activate()
$loginSheet.Cells.Select()

# ---------------------------------------------------------------------
# 3. Add the new "MyAccount" worksheet after "Login".
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$myAccount = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$myAccount.Name = "MyAccount"

# Header row
$myAccount.Range("A1").Value = "TestName"
$myAccount.Range("B1").Value = "Sanity"
$myAccount.Range("C1").Value = "Regression"
$myAccount.Range("D1").Value = "Functional"
$myAccount.Range("E1").Value = "Email"
$myAccount.Range("F1").Value = "Password"
$myAccount.Range("G1").Value = "NewPassword"
$myAccount.Range("H1").Value = "Expected_Error"
$myAccount.Range("A1:H1").Font.Bold = $true

# Row 2 - same old/new password negative flow
$myAccount.Range("A2").Value = "TC_appchangepassword_samenegativeflow_1"
$myAccount.Range("B2").Value = "NO"
$myAccount.Range("C2").Value = "YES"
$myAccount.Range("D2").Value = "NO"
$myAccount.Range("E2").Value = "jqetgFWHup@gmail.com"
$myAccount.Range("F2").Value = "Sgidfn345@"
$myAccount.Range("H2").Value = "Old password and new password should not be same"

# Row 3 - missing upper case char
$myAccount.Range("A3").Value = "TC_appchangepassword_negativeflow_1"
$myAccount.Range("B3").Value = "NO"
$myAccount.Range("C3").Value = "YES"
$myAccount.Range("D3").Value = "NO"
$myAccount.Range("E3").Value = "jqetgFWHup@gmail.com"
$myAccount.Range("F3").Value = "Sgidfn345@"
$myAccount.Range("G3").Value = "asdfadsfas"
$myAccount.Range("H3").Value = "At least one upper case char"

# Row 4 - missing lower case char (new password is purely numeric -> stored as text)
$myAccount.Range("A4").Value = "TC_appchangepassword_negativeflow_2"
$myAccount.Range("B4").Value = "NO"
$myAccount.Range("C4").Value = "YES"
$myAccount.Range("D4").Value = "NO"
$myAccount.Range("E4").Value = "jqetgFWHup@gmail.com"
$myAccount.Range("F4").Value = "Sgidfn345@"
$myAccount.Range("G4").Value = "'23443212"
$myAccount.Range("H4").Value = "At least one lower case char"

# Row 5 - missing non alpha numeric char
$myAccount.Range("A5").Value = "TC_appchangepassword_negativeflow_3"
$myAccount.Range("B5").Value = "NO"
$myAccount.Range("C5").Value = "YES"
$myAccount.Range("D5").Value = "NO"
$myAccount.Range("E5").Value = "jqetgFWHup@gmail.com"
$myAccount.Range("F5").Value = "Sgidfn345@"
$myAccount.Range("G5").Value = "Nsdgs234"
$myAccount.Range("H5").Value = "At least one non alpha numeric char (@,!,#, etc)"

# Row 6 - missing number digit
$myAccount.Range("A6").Value = "TC_appchangepassword_negativeflow_4"
$myAccount.Range("B6").Value = "NO"
$myAccount.Range("C6").Value = "YES"
$myAccount.Range("D6").Value = "NO"
$myAccount.Range("E6").Value = "jqetgFWHup@gmail.com"
$myAccount.Range("F6").Value = "Sgidfn345@"
$myAccount.Range("G6").Value = "Jsdfsdg@"
$myAccount.Range("H6").Value = "At least one number digit"

# Row 7 - too short
$myAccount.Range("A7").Value = "TC_appchangepassword_negativeflow_5"
$myAccount.Range("B7").Value = "NO"
$myAccount.Range("C7").Value = "YES"
$myAccount.Range("D7").Value = "NO"
$myAccount.Range("E7").Value = "jqetgFWHup@gmail.com"
$myAccount.Range("F7").Value = "Sgidfn345@"
$myAccount.Range("G7").Value = "asd"
$myAccount.Range("H7").Value = "Should be at least 8 chars"

# Row 8 - wrong old password
$myAccount.Range("A8").Value = "TC_appchangepassword_erroroldpasswordflow_1"
$myAccount.Range("B8").Value = "NO"
$myAccount.Range("C8").Value = "YES"
$myAccount.Range("D8").Value = "NO"
$myAccount.Range("E8").Value = "jqetgFWHup@gmail.com"
$myAccount.Range("F8").Value = "Sgidfn345@"
$myAccount.Range("G8").Value = "Sgidfe345@"
$myAccount.Range("H8").Value = "Incorrect old password, please try again with the correct password"

# ---------------------------------------------------------------------
# 4. Hyperlinks on the Email / Password / NewPassword columns (mailto:),
#    matching the pattern already used on the other sheets.
# ---------------------------------------------------------------------
$myAccount.Hyperlinks.Add($myAccount.Range("F2"), "mailto:Sgidfn345@")
$myAccount.Hyperlinks.Add($myAccount.Range("E2"), "mailto:jqetgFWHup@gmail.com")
$myAccount.Hyperlinks.Add($myAccount.Range("E3"), "mailto:jqetgFWHup@gmail.com")
$myAccount.Hyperlinks.Add($myAccount.Range("F3"), "mailto:Sgidfn345@")
$myAccount.Hyperlinks.Add($myAccount.Range("G6"), "mailto:Jsdfsdg@")
$myAccount.Hyperlinks.Add($myAccount.Range("E4"), "mailto:jqetgFWHup@gmail.com")
$myAccount.Hyperlinks.Add($myAccount.Range("E5"), "mailto:jqetgFWHup@gmail.com")
$myAccount.Hyperlinks.Add($myAccount.Range("E6"), "mailto:jqetgFWHup@gmail.com")
$myAccount.Hyperlinks.Add($myAccount.Range("F4"), "mailto:Sgidfn345@")
$myAccount.Hyperlinks.Add($myAccount.Range("F5"), "mailto:Sgidfn345@")
$myAccount.Hyperlinks.Add($myAccount.Range("F6"), "mailto:Sgidfn345@")
$myAccount.Hyperlinks.Add($myAccount.Range("E7"), "mailto:jqetgFWHup@gmail.com")
$myAccount.Hyperlinks.Add($myAccount.Range("F7"), "mailto:Sgidfn345@")
$myAccount.Hyperlinks.Add($myAccount.Range("E8"), "mailto:jqetgFWHup@gmail.com")
$myAccount.Hyperlinks.Add($myAccount.Range("F8"), "mailto:Sgidfn345@")
$myAccount.Hyperlinks.Add($myAccount.Range("G8"), "mailto:Sgidfe345@")

# Hyperlinks.Add() re-stamps its own font style; reapply the plain
# "Hyperlink" cell style everywhere it belongs (incl. the styled-but-
# empty G2 cell, which never got a real hyperlink).
$myAccount.Range("E2:F8").Style = "Hyperlink"
$myAccount.Range("G6").Style = "Hyperlink"
$myAccount.Range("G8").Style = "Hyperlink"
$myAccount.Range("G2").Style = "Hyperlink"

$myAccount.PageSetup.Orientation = 1

$myAccount.Range("C2:C8").Select()

# ---------------------------------------------------------------------
# 5. Column widths to roughly match the other data-entry sheets.
# ---------------------------------------------------------------------
$myAccount.Columns("A").ColumnWidth = 39.44140625
$myAccount.Columns("E").ColumnWidth = 29.5546875
$myAccount.Columns("F").ColumnWidth = 20.109375
$myAccount.Columns("G").ColumnWidth = 12.77734375
$myAccount.Columns("H").ColumnWidth = 47.5546875

# ---------------------------------------------------------------------
# 6. Finish back on TestBase, mirroring the final saved selection/tab.
# ---------------------------------------------------------------------
$testBase.Activate()
$testBase.Range("C4").Select()
